# Auto-generated edit script applying the Kujata_Profits market-data refresh diff.
# For each affected cell: set new value, or clear the cell if the diff removed it.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 310
$ws.Range("I11").Value = 310
$ws.Range("K11").Value = 310
$ws.Range("M11").Value = -170
$ws.Range("H15").Value = 3990.3618
$ws.Range("I15").Value = 3990.3618
$ws.Range("K15").Value = 11971.0854
$ws.Range("M15").Value = -11802.0854
$ws.Range("H98").Value = 2443.6667
$ws.Range("I98").Value = 2036.1428
$ws.Range("K98").Value = 2036.1428
$ws.Range("M98").Value = -538.1428000000001
$ws.Range("H122").Value = 2443.6667
$ws.Range("I122").Value = 2036.1428
$ws.Range("K122").Value = 6108.428400000001
$ws.Range("M122").Value = -3658.428400000001
$ws.Range("H132").Value = 7942559.5
$ws.Range("I132").Value = 11907146
$ws.Range("K132").Value = 35721438
$ws.Range("M132").Value = -35718908
$ws.Range("H135").Value = 288.83334
$ws.Range("I135").Value = 246.6
$ws.Range("J135").Value = 500
$ws.Range("K135").Value = 2219.4
$ws.Range("L135").Value = 4500
$ws.Range("M135").Value = 315.5999999999999
$ws.Range("N135").Value = -9570
$ws.Range("H137").Value = 1500
$ws.Range("I137").Value = 1666.6666
$ws.Range("J137").Value = 1250
$ws.Range("K137").Value = 4999.9998
$ws.Range("L137").Value = 3750
$ws.Range("M137").Value = -2449.9998
$ws.Range("N137").Value = -8850
$ws.Range("H138").Value = 751561.6
$ws.Range("I138").Value = 1445.3077
$ws.Range("J138").Value = 968261.9
$ws.Range("K138").Value = 4335.9231
$ws.Range("L138").Value = 2904785.7
$ws.Range("M138").Value = 804.0769
$ws.Range("N138").Value = -2915065.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3171.7
$ws.Range("I132").Value = 2795.8
$ws.Range("K132").Value = 8387.400000000001
$ws.Range("M132").Value = -5857.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 436.5
$ws.Range("I22").Value = 465
$ws.Range("J22").Value = 351
$ws.Range("K22").Value = 465
$ws.Range("L22").Value = 351
$ws.Range("M22").Value = -115
$ws.Range("N22").Value = -1051
$ws.Range("H31").Value = 1967.875
$ws.Range("I31").Value = 1967.875
$ws.Range("K31").Value = 1967.875
$ws.Range("M31").Value = -1672.875
$ws.Range("H34").Value = 1967.875
$ws.Range("I34").Value = 1967.875
$ws.Range("K34").Value = 1967.875
$ws.Range("M34").Value = -1765.875
$ws.Range("H50").Value = 21000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 21000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 21000
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -22250
$ws.Range("H58").Value = 693.41174
$ws.Range("I58").Value = 664
$ws.Range("J58").Value = 789
$ws.Range("K58").Value = 664
$ws.Range("L58").Value = 789
$ws.Range("M58").Value = -461
$ws.Range("N58").Value = -1195
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H81").Value = 16500
$ws.Range("J81").Value = 16500
$ws.Range("L81").Value = 16500
$ws.Range("N81").Value = -18496
$ws.Range("H84").Value = 16500
$ws.Range("J84").Value = 16500
$ws.Range("L84").Value = 49500
$ws.Range("N84").Value = -59484
$ws.Range("H86").Value = 3936193
$ws.Range("I86").Value = 8336033
$ws.Range("J86").Value = 25224
$ws.Range("K86").Value = 8336033
$ws.Range("L86").Value = 25224
$ws.Range("M86").Value = -8334910
$ws.Range("N86").Value = -27470
$ws.Range("H87").Value = 20000
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22372
$ws.Range("H88").Value = 17500
$ws.Range("J88").Value = 17500
$ws.Range("L88").Value = 17500
$ws.Range("N88").Value = -18312
$ws.Range("H89").Value = 3936193
$ws.Range("I89").Value = 8336033
$ws.Range("J89").Value = 25224
$ws.Range("K89").Value = 41680165
$ws.Range("L89").Value = 126120
$ws.Range("M89").Value = -41674549
$ws.Range("N89").Value = -137352
$ws.Range("H90").Value = 20000
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -71856
$ws.Range("H91").Value = 17500
$ws.Range("J91").Value = 17500
$ws.Range("L91").Value = 17500
$ws.Range("N91").Value = -20308
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("H105").Value = 799.75
$ws.Range("I105").Value = 799.75
$ws.Range("K105").Value = 799.75
$ws.Range("M105").Value = 947.25
$ws.Range("H107").Value = 645
$ws.Range("I107").Value = 414
$ws.Range("K107").Value = 414
$ws.Range("M107").Value = 1506
$ws.Range("H132").Value = 8864.471
$ws.Range("I132").Value = 13448.333
$ws.Range("J132").Value = 3707.625
$ws.Range("K132").Value = 40344.999
$ws.Range("L132").Value = 11122.875
$ws.Range("M132").Value = -37814.999
$ws.Range("N132").Value = -16182.875
$ws.Range("H134").Value = 20835406
$ws.Range("I134").Value = 27779750
$ws.Range("K134").Value = 83339250
$ws.Range("M134").Value = -83336715
$ws.Range("H136").Value = 693.41174
$ws.Range("I136").Value = 664
$ws.Range("J136").Value = 789
$ws.Range("K136").Value = 1992
$ws.Range("L136").Value = 2367
$ws.Range("M136").Value = 558
$ws.Range("N136").Value = -7467

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 780.2
$ws.Range("I122").Value = 702.7143
$ws.Range("J122").Value = 848
$ws.Range("K122").Value = 6324.428699999999
$ws.Range("L122").Value = 7632
$ws.Range("M122").Value = -3874.428699999999
$ws.Range("N122").Value = -12532
$ws.Range("H131").Value = 20001376
$ws.Range("J131").Value = 1591.3658
$ws.Range("L131").Value = 4774.097400000001
$ws.Range("N131").Value = -14854.0974
$ws.Range("H140").Value = 40617.242
$ws.Range("I140").Value = 48454.332
$ws.Range("J140").Value = 2999.2
$ws.Range("K140").Value = 145362.996
$ws.Range("L140").Value = 8997.599999999999
$ws.Range("M140").Value = -140182.996
$ws.Range("N140").Value = -19357.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H55").Value = 3033
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H102").Value = 2626.9092
$ws.Range("I102").Value = 2889.5
$ws.Range("J102").Value = 1926.6666
$ws.Range("K102").Value = 2889.5
$ws.Range("L102").Value = 1926.6666
$ws.Range("M102").Value = -1267.5
$ws.Range("N102").Value = -5170.6666
$ws.Range("H132").Value = 2352.2727
$ws.Range("I132").Value = 2199.4
$ws.Range("J132").Value = 2679.8572
$ws.Range("K132").Value = 6598.200000000001
$ws.Range("L132").Value = 8039.571599999999
$ws.Range("M132").Value = -4068.200000000001
$ws.Range("N132").Value = -13099.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1081.4117
$ws.Range("I22").Value = 1216.3636
$ws.Range("J22").Value = 834
$ws.Range("K22").Value = 1216.3636
$ws.Range("L22").Value = 834
$ws.Range("M22").Value = -921.3635999999999
$ws.Range("N22").Value = -1424
$ws.Range("H27").Value = 1081.4117
$ws.Range("I27").Value = 1216.3636
$ws.Range("J27").Value = 834
$ws.Range("K27").Value = 1216.3636
$ws.Range("L27").Value = 834
$ws.Range("M27").Value = -1109.3636
$ws.Range("N27").Value = -1048
$ws.Range("H61").Value = 1556.8334
$ws.Range("I61").Value = 1394.6364
$ws.Range("K61").Value = 1394.6364
$ws.Range("M61").Value = -1192.6364
$ws.Range("H68").Value = 1413.5714
$ws.Range("I68").Value = 1220.4
$ws.Range("J68").Value = 1896.5
$ws.Range("K68").Value = 1220.4
$ws.Range("L68").Value = 1896.5
$ws.Range("M68").Value = -471.4000000000001
$ws.Range("N68").Value = -3394.5
$ws.Range("H71").Value = 1413.5714
$ws.Range("I71").Value = 1220.4
$ws.Range("J71").Value = 1896.5
$ws.Range("K71").Value = 6102
$ws.Range("L71").Value = 9482.5
$ws.Range("M71").Value = -2358
$ws.Range("N71").Value = -16970.5
$ws.Range("H87").Value = 13000
$ws.Range("I87").Value = 13000
$ws.Range("K87").Value = 13000
$ws.Range("M87").Value = -11877
$ws.Range("H90").Value = 13000
$ws.Range("I90").Value = 13000
$ws.Range("K90").Value = 39000
$ws.Range("M90").Value = -33384
$ws.Range("H111").Value = 66666.664
$ws.Range("J111").Value = 66666.664
$ws.Range("L111").Value = 66666.664
$ws.Range("N111").Value = -74846.664
$ws.Range("H113").Value = 1556.8334
$ws.Range("I113").Value = 1394.6364
$ws.Range("K113").Value = 1394.6364
$ws.Range("M113").Value = 775.3635999999999
$ws.Range("H132").Value = 103570.1
$ws.Range("I132").Value = 3966.3333
$ws.Range("J132").Value = 146257.42
$ws.Range("K132").Value = 11898.9999
$ws.Range("L132").Value = 438772.26
$ws.Range("M132").Value = -9368.999899999999
$ws.Range("N132").Value = -443832.26

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 487.78946
$ws.Range("I113").Value = 373.66666
$ws.Range("K113").Value = 1120.99998
$ws.Range("M113").Value = 1049.00002
$ws.Range("H132").Value = 4600.0527
$ws.Range("I132").Value = 4257.4287
$ws.Range("K132").Value = 12772.2861
$ws.Range("M132").Value = -10242.2861
